# The "Recorded By" column (G) lists the users/processes that recorded a
# session. Previously "System" could appear anywhere in the comma
# separated list; it should always be listed first.
#
# This walks every populated row of the "Recorded By" column and, if the
# token "System" is present but not already the first entry, moves it to
# the front of the list while preserving the relative order of the other
# entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") { continue }

    $parts = @($text -split ', ')
    $idx = [Array]::IndexOf($parts, "System")

    if ($idx -gt 0) {
        $rest = @()
        for ($i = 0; $i -lt $parts.Count; $i++) {
            if ($i -ne $idx) {
                $rest += $parts[$i]
            }
        }
        $newParts = @("System") + $rest
        $cell.Value = [string]::Join(", ", $newParts)
    }
}
